$d = $word.ActiveDocument

$replacements = @(
    @("490÷4=", "400÷2="),
    @("191÷6=", "947÷7="),
    @("769÷8=", "859÷2="),
    @("183÷6=", "682÷9="),
    @("944÷8=", "305÷5="),
    @("295÷2=", "804÷8="),
    @("905÷4=", "418÷4="),
    @("496÷4=", "364÷7="),
    @("242÷4=", "285÷3="),
    @("741÷6=", "899÷9="),
    @("276÷4=", "224÷8="),
    @("525÷8=", "545÷9="),
    @("254÷9=", "732÷6="),
    @("997÷5=", "707÷8="),
    @("147÷8=", "681÷6="),
    @("460÷8=", "503÷4="),
    @("821÷5=", "420÷5="),
    @("839÷3=", "937÷5="),
    @("216÷3=", "549÷7="),
    @("647÷3=", "716÷5="),
    @("707÷7=", "161÷2="),
    @("310÷3=", "195÷7="),
    @("341÷8=", "992÷6="),
    @("113÷8=", "979÷6="),
    @("269÷2=", "281÷8=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
